$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 64,6
$data[0,0] = "Identifying the Human Rights Impacts of Palm Oil"
$data[0,1] = "Dr. Maria Schmidt"
$data[0,2] = "A1"
$data[0,3] = 8.4
$data[0,4] = "2025-07-03 03:27:14"
$data[0,5] = 10
$data[1,0] = "Handlungsleitfaden zur Durchführung der Risikoanalyse für Kakao produzierende Länder"
$data[1,1] = "Prof. Dr. James Allen"
$data[1,2] = "D35"
$data[1,3] = 6.7
$data[1,4] = "2025-07-03 03:27:14"
$data[1,5] = 7
$data[2,0] = "UNEP FI Human Rights Guidance Tool for the Financial Sector `"Infrastructure`""
$data[2,1] = "Abdallah Reyati"
$data[2,2] = "H50"
$data[2,3] = 8
$data[2,4] = "2025-07-03 03:27:14"
$data[2,5] = 4
$data[3,0] = "Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries"
$data[3,1] = "Lisa Müller"
$data[3,2] = "C15"
$data[3,3] = 5.7
$data[3,4] = "2025-07-03 03:27:14"
$data[3,5] = 10
$data[4,0] = "Praxislotse Wirtschaft & Menschenrechte"
$data[4,1] = "Dr. Maria Schmidt"
$data[4,2] = "D35"
$data[4,3] = 8.9
$data[4,4] = "2025-07-03 03:27:14"
$data[4,5] = 7
$data[5,0] = "UNEP FI Human Rights Guidance Tool for the Financial Sector `"Oil and Gas`""
$data[5,1] = "Prof. Dr. James Allen"
$data[5,2] = "C13"
$data[5,3] = 6.4
$data[5,4] = "2025-07-03 03:27:14"
$data[5,5] = 4
$data[6,0] = "UNEP FI Human Rights Guidance Tool for the Financial Sector `"Utilities and Waste Management`""
$data[6,1] = "Abdallah Reyati"
$data[6,2] = "E38"
$data[6,3] = 8.699999999999999
$data[6,4] = "2025-07-03 03:27:14"
$data[6,5] = 10
$data[7,0] = "Strengthening Protection Against Trafficking in Persons in Federal and Corporate Supply Chains"
$data[7,1] = "Lisa Müller"
$data[7,2] = "A1"
$data[7,3] = 5.7
$data[7,4] = "2025-07-03 03:27:14"
$data[7,5] = 7
$data[8,0] = "Human Rights Risks in Mining, A Baseline Study"
$data[8,1] = "Dr. Maria Schmidt"
$data[8,2] = "A1"
$data[8,3] = 8.6
$data[8,4] = "2025-07-03 03:27:14"
$data[8,5] = 4
$data[9,0] = "2022 List Of Goods Produced By Child Labor Or Forced Labor"
$data[9,1] = "Prof. Dr. James Allen"
$data[9,2] = "C21"
$data[9,3] = 6.3
$data[9,4] = "2025-07-03 03:27:14"
$data[9,5] = 10
$data[10,0] = "Die Achtung von Menschenrechten entlang globaler Wertschöpfungskette"
$data[10,1] = "Abdallah Reyati"
$data[10,2] = "D35"
$data[10,3] = 8.300000000000001
$data[10,4] = "2025-07-03 03:27:14"
$data[10,5] = 7
$data[11,0] = "Global Estimates of Modern Slavery, Forced Labour and Forced Marriage"
$data[11,1] = "Lisa Müller"
$data[11,2] = "A1"
$data[11,3] = 5
$data[11,4] = "2025-07-03 03:27:14"
$data[11,5] = 4
$data[12,0] = "Low Prices Drive Natural Rubber Producers Into Poverty"
$data[12,1] = "Dr. Maria Schmidt"
$data[12,2] = "D35"
$data[12,3] = 8.4
$data[12,4] = "2025-07-03 03:27:14"
$data[12,5] = 10
$data[13,0] = "Global Mercury Assessment 2018"
$data[13,1] = "Prof. Dr. James Allen"
$data[13,2] = "D35"
$data[13,3] = 6.7
$data[13,4] = "2025-07-03 03:27:14"
$data[13,5] = 7
$data[14,0] = "Sustainability and Circularity in the Textile Value Chain - A Global Roadmap"
$data[14,1] = "Abdallah Reyati"
$data[14,2] = "C13"
$data[14,3] = 8
$data[14,4] = "2025-07-03 03:27:14"
$data[14,5] = 4
$data[15,0] = "Forced and Child Labour in the Cotton Industry"
$data[15,1] = "Lisa Müller"
$data[15,2] = "A1"
$data[15,3] = 6.1
$data[15,4] = "2025-07-03 03:27:14"
$data[15,5] = 10
$data[16,0] = "Forced, child and trafficked labour in the cocoa industry"
$data[16,1] = "Dr. Maria Schmidt"
$data[16,2] = "A1"
$data[16,3] = 8.4
$data[16,4] = "2025-07-03 03:27:14"
$data[16,5] = 7
$data[17,0] = "Ermittlung von potentiell POP-haltigen Abfällen und Recyclingstoffen - Ableitung von Grenzwerten"
$data[17,1] = "Prof. Dr. James Allen"
$data[17,2] = "D35"
$data[17,3] = 6.4
$data[17,4] = "2025-07-03 03:27:14"
$data[17,5] = 4
$data[18,0] = "Kurzzeitige Chlorparafine (SCCP) Vorkommen, Verwendung und Rechtssetzung zu kurzketigen Chlorparaffinen in Produkten und Abfällen"
$data[18,1] = "Abdallah Reyati"
$data[18,2] = "A1"
$data[18,3] = 7.9
$data[18,4] = "2025-07-03 03:27:14"
$data[18,5] = 10
$data[19,0] = "Hexabromocyclododecane"
$data[19,1] = "Lisa Müller"
$data[19,2] = "D35"
$data[19,3] = 5.3
$data[19,4] = "2025-07-03 03:27:14"
$data[19,5] = 7
$data[20,0] = "Fairtrade Risk Map"
$data[20,1] = "Dr. Maria Schmidt"
$data[20,2] = "D35"
$data[20,3] = 8.1
$data[20,4] = "2025-07-03 03:27:14"
$data[20,5] = 4
$data[21,0] = "Typical Wastes Generated By Industry Sector"
$data[21,1] = "Prof. Dr. James Allen"
$data[21,2] = "C21"
$data[21,3] = 7.1
$data[21,4] = "2025-07-03 03:27:14"
$data[21,5] = 10
$data[22,0] = "Cocoa Barometer 2020"
$data[22,1] = "Abdallah Reyati"
$data[22,2] = "A1"
$data[22,3] = 7.9
$data[22,4] = "2025-07-03 03:27:14"
$data[22,5] = 7
$data[23,0] = "Handlungsleitfaden zur Durchführung der Risikoanalysen für Kakao produzierende Länder"
$data[23,1] = "Lisa Müller"
$data[23,2] = "D35"
$data[23,3] = 5
$data[23,4] = "2025-07-03 03:27:14"
$data[23,5] = 4
$data[24,0] = "Human Rights Toolkit for Financial Institutions"
$data[24,1] = "Dr. Maria Schmidt"
$data[24,2] = "C21"
$data[24,3] = 9.300000000000001
$data[24,4] = "2025-07-03 03:27:14"
$data[24,5] = 10
$data[25,0] = "Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries"
$data[25,1] = "Prof. Dr. James Allen"
$data[25,2] = "C15"
$data[25,3] = 5.9
$data[25,4] = "2025-07-03 03:27:14"
$data[25,5] = 7
$data[26,0] = "Global Dialogue Forum on Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries"
$data[26,1] = "Abdallah Reyati"
$data[26,2] = "C15"
$data[26,3] = 7.1
$data[26,4] = "2025-07-03 03:27:14"
$data[26,5] = 4
$data[27,0] = "Praxislotse Wirtschaft und Menschenrechte"
$data[27,1] = "Lisa Müller"
$data[27,2] = "D35"
$data[27,3] = 6.6
$data[27,4] = "2025-07-03 03:27:14"
$data[27,5] = 10
$data[28,0] = "Verité Commodity Atlas"
$data[28,1] = "Dr. Maria Schmidt"
$data[28,2] = "D35"
$data[28,3] = 8.4
$data[28,4] = "2025-07-03 03:27:14"
$data[28,5] = 7
$data[29,0] = "Business & Human Rights Navigator"
$data[29,1] = "Prof. Dr. James Allen"
$data[29,2] = "D35"
$data[29,3] = 6.4
$data[29,4] = "2025-07-03 03:27:14"
$data[29,5] = 4
$data[30,0] = "Losing Ground, The Human Rights Impacts of Oil Palm Plantation Expansion in Indonesia"
$data[30,1] = "Abdallah Reyati"
$data[30,2] = "A1"
$data[30,3] = 7.9
$data[30,4] = "2025-07-03 03:27:14"
$data[30,5] = 10
$data[31,0] = "When We Lost the Forest, We Lost Everything: Oil Palm Plantations and Rights Violations in Indonesia"
$data[31,1] = "Lisa Müller"
$data[31,2] = "A1"
$data[31,3] = 5.3
$data[31,4] = "2025-07-03 03:27:14"
$data[31,5] = 7
$data[32,0] = "Identifying the Human Rights Impacts of Palm Oil"
$data[32,1] = "Dr. Maria Schmidt"
$data[32,2] = "A1"
$data[32,3] = 8.4
$data[32,4] = "2025-07-03 03:56:32"
$data[32,5] = 10
$data[33,0] = "Handlungsleitfaden zur Durchführung der Risikoanalyse für Kakao produzierende Länder"
$data[33,1] = "Prof. Dr. James Allen"
$data[33,2] = "D35"
$data[33,3] = 6.7
$data[33,4] = "2025-07-03 03:56:32"
$data[33,5] = 7
$data[34,0] = "UNEP FI Human Rights Guidance Tool for the Financial Sector `"Infrastructure`""
$data[34,1] = "Abdallah Reyati"
$data[34,2] = "H50"
$data[34,3] = 8
$data[34,4] = "2025-07-03 03:56:32"
$data[34,5] = 4
$data[35,0] = "Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries"
$data[35,1] = "Lisa Müller"
$data[35,2] = "C15"
$data[35,3] = 5.7
$data[35,4] = "2025-07-03 03:56:32"
$data[35,5] = 10
$data[36,0] = "Praxislotse Wirtschaft & Menschenrechte"
$data[36,1] = "Dr. Maria Schmidt"
$data[36,2] = "D35"
$data[36,3] = 8.9
$data[36,4] = "2025-07-03 03:56:32"
$data[36,5] = 7
$data[37,0] = "UNEP FI Human Rights Guidance Tool for the Financial Sector `"Oil and Gas`""
$data[37,1] = "Prof. Dr. James Allen"
$data[37,2] = "C13"
$data[37,3] = 6.4
$data[37,4] = "2025-07-03 03:56:32"
$data[37,5] = 4
$data[38,0] = "UNEP FI Human Rights Guidance Tool for the Financial Sector `"Utilities and Waste Management`""
$data[38,1] = "Abdallah Reyati"
$data[38,2] = "E38"
$data[38,3] = 8.699999999999999
$data[38,4] = "2025-07-03 03:56:32"
$data[38,5] = 10
$data[39,0] = "Strengthening Protection Against Trafficking in Persons in Federal and Corporate Supply Chains"
$data[39,1] = "Lisa Müller"
$data[39,2] = "A1"
$data[39,3] = 5.7
$data[39,4] = "2025-07-03 03:56:32"
$data[39,5] = 7
$data[40,0] = "Human Rights Risks in Mining, A Baseline Study"
$data[40,1] = "Dr. Maria Schmidt"
$data[40,2] = "A1"
$data[40,3] = 8.6
$data[40,4] = "2025-07-03 03:56:32"
$data[40,5] = 4
$data[41,0] = "2022 List Of Goods Produced By Child Labor Or Forced Labor"
$data[41,1] = "Prof. Dr. James Allen"
$data[41,2] = "C21"
$data[41,3] = 6.3
$data[41,4] = "2025-07-03 03:56:32"
$data[41,5] = 10
$data[42,0] = "Die Achtung von Menschenrechten entlang globaler Wertschöpfungskette"
$data[42,1] = "Abdallah Reyati"
$data[42,2] = "D35"
$data[42,3] = 8.300000000000001
$data[42,4] = "2025-07-03 03:56:32"
$data[42,5] = 7
$data[43,0] = "Global Estimates of Modern Slavery, Forced Labour and Forced Marriage"
$data[43,1] = "Lisa Müller"
$data[43,2] = "A1"
$data[43,3] = 5
$data[43,4] = "2025-07-03 03:56:32"
$data[43,5] = 4
$data[44,0] = "Low Prices Drive Natural Rubber Producers Into Poverty"
$data[44,1] = "Dr. Maria Schmidt"
$data[44,2] = "D35"
$data[44,3] = 8.4
$data[44,4] = "2025-07-03 03:56:32"
$data[44,5] = 10
$data[45,0] = "Global Mercury Assessment 2018"
$data[45,1] = "Prof. Dr. James Allen"
$data[45,2] = "D35"
$data[45,3] = 6.7
$data[45,4] = "2025-07-03 03:56:32"
$data[45,5] = 7
$data[46,0] = "Sustainability and Circularity in the Textile Value Chain - A Global Roadmap"
$data[46,1] = "Abdallah Reyati"
$data[46,2] = "C13"
$data[46,3] = 8
$data[46,4] = "2025-07-03 03:56:32"
$data[46,5] = 4
$data[47,0] = "Forced and Child Labour in the Cotton Industry"
$data[47,1] = "Lisa Müller"
$data[47,2] = "A1"
$data[47,3] = 6.1
$data[47,4] = "2025-07-03 03:56:32"
$data[47,5] = 10
$data[48,0] = "Forced, child and trafficked labour in the cocoa industry"
$data[48,1] = "Dr. Maria Schmidt"
$data[48,2] = "A1"
$data[48,3] = 8.4
$data[48,4] = "2025-07-03 03:56:32"
$data[48,5] = 7
$data[49,0] = "Ermittlung von potentiell POP-haltigen Abfällen und Recyclingstoffen - Ableitung von Grenzwerten"
$data[49,1] = "Prof. Dr. James Allen"
$data[49,2] = "D35"
$data[49,3] = 6.4
$data[49,4] = "2025-07-03 03:56:32"
$data[49,5] = 4
$data[50,0] = "Kurzzeitige Chlorparafine (SCCP) Vorkommen, Verwendung und Rechtssetzung zu kurzketigen Chlorparaffinen in Produkten und Abfällen"
$data[50,1] = "Abdallah Reyati"
$data[50,2] = "A1"
$data[50,3] = 7.9
$data[50,4] = "2025-07-03 03:56:32"
$data[50,5] = 10
$data[51,0] = "Hexabromocyclododecane"
$data[51,1] = "Lisa Müller"
$data[51,2] = "D35"
$data[51,3] = 5.3
$data[51,4] = "2025-07-03 03:56:32"
$data[51,5] = 7
$data[52,0] = "Fairtrade Risk Map"
$data[52,1] = "Dr. Maria Schmidt"
$data[52,2] = "D35"
$data[52,3] = 8.1
$data[52,4] = "2025-07-03 03:56:32"
$data[52,5] = 4
$data[53,0] = "Typical Wastes Generated By Industry Sector"
$data[53,1] = "Prof. Dr. James Allen"
$data[53,2] = "C21"
$data[53,3] = 7.1
$data[53,4] = "2025-07-03 03:56:32"
$data[53,5] = 10
$data[54,0] = "Cocoa Barometer 2020"
$data[54,1] = "Abdallah Reyati"
$data[54,2] = "A1"
$data[54,3] = 7.9
$data[54,4] = "2025-07-03 03:56:32"
$data[54,5] = 7
$data[55,0] = "Handlungsleitfaden zur Durchführung der Risikoanalysen für Kakao produzierende Länder"
$data[55,1] = "Lisa Müller"
$data[55,2] = "D35"
$data[55,3] = 5
$data[55,4] = "2025-07-03 03:56:32"
$data[55,5] = 4
$data[56,0] = "Human Rights Toolkit for Financial Institutions"
$data[56,1] = "Dr. Maria Schmidt"
$data[56,2] = "C21"
$data[56,3] = 9.300000000000001
$data[56,4] = "2025-07-03 03:56:32"
$data[56,5] = 10
$data[57,0] = "Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries"
$data[57,1] = "Prof. Dr. James Allen"
$data[57,2] = "C15"
$data[57,3] = 5.9
$data[57,4] = "2025-07-03 03:56:32"
$data[57,5] = 7
$data[58,0] = "Global Dialogue Forum on Wages and Working Hours in the Textiles, Clothing, Leather and Footwear Industries"
$data[58,1] = "Abdallah Reyati"
$data[58,2] = "C15"
$data[58,3] = 7.1
$data[58,4] = "2025-07-03 03:56:32"
$data[58,5] = 4
$data[59,0] = "Praxislotse Wirtschaft und Menschenrechte"
$data[59,1] = "Lisa Müller"
$data[59,2] = "D35"
$data[59,3] = 6.6
$data[59,4] = "2025-07-03 03:56:32"
$data[59,5] = 10
$data[60,0] = "Verité Commodity Atlas"
$data[60,1] = "Dr. Maria Schmidt"
$data[60,2] = "D35"
$data[60,3] = 8.4
$data[60,4] = "2025-07-03 03:56:32"
$data[60,5] = 7
$data[61,0] = "Business & Human Rights Navigator"
$data[61,1] = "Prof. Dr. James Allen"
$data[61,2] = "D35"
$data[61,3] = 6.4
$data[61,4] = "2025-07-03 03:56:32"
$data[61,5] = 4
$data[62,0] = "Losing Ground, The Human Rights Impacts of Oil Palm Plantation Expansion in Indonesia"
$data[62,1] = "Abdallah Reyati"
$data[62,2] = "A1"
$data[62,3] = 7.9
$data[62,4] = "2025-07-03 03:56:32"
$data[62,5] = 10
$data[63,0] = "When We Lost the Forest, We Lost Everything: Oil Palm Plantations and Rights Violations in Indonesia"
$data[63,1] = "Lisa Müller"
$data[63,2] = "A1"
$data[63,3] = 5.3
$data[63,4] = "2025-07-03 03:56:32"
$data[63,5] = 7

$ws.Range("A132:F195").Value = $data
